# Clear the "Font Weight" values in column I for rows 2-13 (the style guide
# rows) while keeping their existing cell style/formatting intact, then
# select the cleared range as the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Style Guide")

$range = $ws.Range("I2:I13")
$range.ClearContents()

$ws.Range("I2:I13").Select()
